# "excel validation bug fixes"
#
# 1. C2 10+10 should equal 21 now (was 20) - updated input value.
# 2. D2 was an (empty) text/shared-string placeholder with a highlighted
#    style; it becomes a plain numeric 20 with no special formatting.
# 3. D5 / D6 lose their manual highlight style (back to the default style).
# 4. A new row (row 7) is appended: "AS" / 1.
# 5. The conditional-formatting "sum check" rules are reworked: they now
#    guard with ISNUMBER() (so blank rows don't light up red) and are
#    extended to the full columns, with the D column split into its own
#    rule set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- cell value fixes ---------------------------------------------------
$ws.Range("C2").Value = 21

$ws.Range("D2").ClearFormats()
$ws.Range("D2").Value = 20

$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()

# --- new row --------------------------------------------------------------
$ws.Range("A7").Value = "AS"
$ws.Range("B7").Value = 1

# --- conditional formatting rework ----------------------------------------
$rngABC = $ws.Range("A1:C1048576")
$rngD = $ws.Range("D1:D1048576")

$green = 10289081
$red = 3424491
$orange = 42495

$fcOk = $rngABC.FormatConditions.Add(2, 0, "=AND(`$C1=`$A1+`$B1, ISNUMBER(`$A1:`$C1))")
$fcOk.Interior.Color = $green

$fcBad = $rngABC.FormatConditions.Add(2, 0, "=AND(`$C1<>`$A1+`$B1, ISNUMBER(`$C1))")
$fcBad.Interior.Color = $red

$fcBadD = $rngD.FormatConditions.Add(2, 0, "=AND(`$C1<>`$A1+`$B1, ISNUMBER(`$C1))")
$fcBadD.Interior.Color = $orange

# --- final selection --------------------------------------------------------
$ws.Range("D7").Select()
